$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 14806
$ws1.Range("F3").Value = 18289
$ws1.Range("F5").Value = 101
$ws1.Range("F14").Value = 85
$ws1.Range("F17").Value = 1387
$ws1.Range("F22").Value = 7573
$ws1.Range("F24").Value = 15
$ws1.Range("F28").Value = 5918
$ws1.Range("F30").Value = 56
$ws1.Range("F34").Value = 5248
$ws1.Range("F36").Value = 36

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("G2").Value = "不可售"

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 14806
$ws4.Range("F3").Value = 18289
$ws4.Range("F5").Value = 101
$ws4.Range("F14").Value = 85
$ws4.Range("F17").Value = 1387
$ws4.Range("F23").Value = 7573
$ws4.Range("F25").Value = 15
$ws4.Range("F30").Value = 5918
$ws4.Range("F32").Value = 56
$ws4.Range("F36").Value = 5248
$ws4.Range("F38").Value = 36
$ws4.Range("G20").Value = "不可售"
